# Auto-generated edit script: apply updated Leve profit calculations
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets ("scheduled runner" price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2639.7932
$ws.Cells.Item(62, 9).Value = 2515.3809
$ws.Cells.Item(62, 10).Value = 2966.375
$ws.Cells.Item(62, 11).Value = 2515.3809
$ws.Cells.Item(62, 12).Value = 2966.375
$ws.Cells.Item(62, 13).Value = -1891.3809
$ws.Cells.Item(62, 14).Value = -4214.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 2639.7932
$ws.Cells.Item(65, 9).Value = 2515.3809
$ws.Cells.Item(65, 10).Value = 2966.375
$ws.Cells.Item(65, 11).Value = 12576.9045
$ws.Cells.Item(65, 12).Value = 14831.875
$ws.Cells.Item(65, 13).Value = -9456.904500000001
$ws.Cells.Item(65, 14).Value = -21071.875

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 992.8261
$ws.Cells.Item(127, 9).Value = 394.0909
$ws.Cells.Item(127, 10).Value = 1541.6666
$ws.Cells.Item(127, 11).Value = 1182.2727
$ws.Cells.Item(127, 12).Value = 4624.9998
$ws.Cells.Item(127, 13).Value = 3777.7273
$ws.Cells.Item(127, 14).Value = -14544.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 4349893
$ws.Cells.Item(137, 9).Value = 7693715.5
$ws.Cells.Item(137, 10).Value = 2924.3
$ws.Cells.Item(137, 11).Value = 23081146.5
$ws.Cells.Item(137, 12).Value = 8772.900000000001
$ws.Cells.Item(137, 13).Value = -23078596.5
$ws.Cells.Item(137, 14).Value = -13872.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 22524.281
$ws.Cells.Item(32, 9).Value = 20488.889
$ws.Cells.Item(32, 10).Value = 33515.4
$ws.Cells.Item(32, 11).Value = 20488.889
$ws.Cells.Item(32, 12).Value = 33515.4
$ws.Cells.Item(32, 13).Value = -20201.889
$ws.Cells.Item(32, 14).Value = -34089.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2910.7334
$ws.Cells.Item(45, 9).Value = 2940.0715
$ws.Cells.Item(45, 11).Value = 2940.0715
$ws.Cells.Item(45, 13).Value = -2563.0715

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 182001310
$ws.Cells.Item(61, 9).Value = 143001390
$ws.Cells.Item(61, 10).Value = 250251200
$ws.Cells.Item(61, 11).Value = 143001390
$ws.Cells.Item(61, 12).Value = 250251200
$ws.Cells.Item(61, 13).Value = -143001178
$ws.Cells.Item(61, 14).Value = -250251624

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 14824341
$ws.Cells.Item(74, 9).Value = 17929358
$ws.Cells.Item(74, 11).Value = 17929358
$ws.Cells.Item(74, 13).Value = -17928484

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 14824341
$ws.Cells.Item(77, 9).Value = 17929358
$ws.Cells.Item(77, 11).Value = 89646790
$ws.Cells.Item(77, 13).Value = -89642422

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value = 55384.617
$ws.Cells.Item(123, 10).Value = 55384.617
$ws.Cells.Item(123, 12).Value = 55384.617
$ws.Cells.Item(123, 14).Value = -65184.617

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 38952.87
$ws.Cells.Item(132, 9).Value = 29448.223
$ws.Cells.Item(132, 11).Value = 88344.66900000001
$ws.Cells.Item(132, 13).Value = -85814.66900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 182001310
$ws.Cells.Item(136, 9).Value = 143001390
$ws.Cells.Item(136, 10).Value = 250251200
$ws.Cells.Item(136, 11).Value = 429004170
$ws.Cells.Item(136, 12).Value = 750753600
$ws.Cells.Item(136, 13).Value = -429001620
$ws.Cells.Item(136, 14).Value = -750758700

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 303.26666
$ws.Cells.Item(22, 9).Value = 196.15384
$ws.Cells.Item(22, 10).Value = 999.5
$ws.Cells.Item(22, 11).Value = 196.15384
$ws.Cells.Item(22, 12).Value = 999.5
$ws.Cells.Item(22, 13).Value = -23.15384
$ws.Cells.Item(22, 14).Value = -1345.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 658.2727
$ws.Cells.Item(94, 9).Value = 440.53333
$ws.Cells.Item(94, 10).Value = 1124.8572
$ws.Cells.Item(94, 11).Value = 440.53333
$ws.Cells.Item(94, 12).Value = 1124.8572
$ws.Cells.Item(94, 13).Value = 10.46667000000002
$ws.Cells.Item(94, 14).Value = -2026.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2022.1666
$ws.Cells.Item(99, 9).Value = 2750
$ws.Cells.Item(99, 10).Value = 1658.25
$ws.Cells.Item(99, 11).Value = 2750
$ws.Cells.Item(99, 12).Value = 1658.25
$ws.Cells.Item(99, 13).Value = -1252
$ws.Cells.Item(99, 14).Value = -4654.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4120.2
$ws.Cells.Item(134, 9).Value = 4700.25
$ws.Cells.Item(134, 11).Value = 14100.75
$ws.Cells.Item(134, 13).Value = -11565.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3056.325
$ws.Cells.Item(31, 10).Value = 7717.636
$ws.Cells.Item(31, 12).Value = 7717.636
$ws.Cells.Item(31, 14).Value = -8307.636

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3056.325
$ws.Cells.Item(34, 10).Value = 7717.636
$ws.Cells.Item(34, 12).Value = 7717.636
$ws.Cells.Item(34, 14).Value = -8121.636

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 41668284
$ws.Cells.Item(58, 9).Value = 41667852
$ws.Cells.Item(58, 10).Value = 41668716
$ws.Cells.Item(58, 11).Value = 41667852
$ws.Cells.Item(58, 12).Value = 41668716
$ws.Cells.Item(58, 13).Value = -41667649
$ws.Cells.Item(58, 14).Value = -41669122

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 44048.625
$ws.Cells.Item(132, 9).Value = 2270.2354
$ws.Cells.Item(132, 10).Value = 145510.42
$ws.Cells.Item(132, 11).Value = 6810.706200000001
$ws.Cells.Item(132, 12).Value = 436531.26
$ws.Cells.Item(132, 13).Value = -4280.706200000001
$ws.Cells.Item(132, 14).Value = -441591.26

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 47153.793
$ws.Cells.Item(134, 9).Value = 2389.0588
$ws.Cells.Item(134, 10).Value = 155868.14
$ws.Cells.Item(134, 11).Value = 7167.176399999999
$ws.Cells.Item(134, 12).Value = 467604.42
$ws.Cells.Item(134, 13).Value = -4632.176399999999
$ws.Cells.Item(134, 14).Value = -472674.42

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 41668284
$ws.Cells.Item(136, 9).Value = 41667852
$ws.Cells.Item(136, 10).Value = 41668716
$ws.Cells.Item(136, 11).Value = 125003556
$ws.Cells.Item(136, 12).Value = 125006148
$ws.Cells.Item(136, 13).Value = -125001006
$ws.Cells.Item(136, 14).Value = -125011248

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 670
$ws.Cells.Item(122, 10).Value = 1496.8
$ws.Cells.Item(122, 12).Value = 13471.2
$ws.Cells.Item(122, 14).Value = -18371.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1303.4681
$ws.Cells.Item(131, 9).Value = 621.875
$ws.Cells.Item(131, 10).Value = 1655.258
$ws.Cells.Item(131, 11).Value = 1865.625
$ws.Cells.Item(131, 12).Value = 4965.774
$ws.Cells.Item(131, 13).Value = 3174.375
$ws.Cells.Item(131, 14).Value = -15045.774

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1973.25
$ws.Cells.Item(102, 9).Value = 1840.6875
$ws.Cells.Item(102, 10).Value = 2503.5
$ws.Cells.Item(102, 11).Value = 1840.6875
$ws.Cells.Item(102, 12).Value = 2503.5
$ws.Cells.Item(102, 13).Value = -218.6875
$ws.Cells.Item(102, 14).Value = -5747.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 12449.5
$ws.Cells.Item(123, 10).Value = 12449.5
$ws.Cells.Item(123, 12).Value = 12449.5
$ws.Cells.Item(123, 14).Value = -17349.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 92941.37
$ws.Cells.Item(132, 9).Value = 112068.22
$ws.Cells.Item(132, 10).Value = 79699.69500000001
$ws.Cells.Item(132, 11).Value = 336204.66
$ws.Cells.Item(132, 12).Value = 239099.085
$ws.Cells.Item(132, 13).Value = -333674.66
$ws.Cells.Item(132, 14).Value = -244159.085

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2751.75
$ws.Cells.Item(122, 9).Value = 2832.8
$ws.Cells.Item(122, 10).Value = 2616.6667
$ws.Cells.Item(122, 11).Value = 8498.400000000001
$ws.Cells.Item(122, 12).Value = 7850.000100000001
$ws.Cells.Item(122, 13).Value = -6048.400000000001
$ws.Cells.Item(122, 14).Value = -12750.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 74932.71000000001
$ws.Cells.Item(132, 9).Value = 1608.5714
$ws.Cells.Item(132, 10).Value = 148256.86
$ws.Cells.Item(132, 11).Value = 4825.7142
$ws.Cells.Item(132, 12).Value = 444770.58
$ws.Cells.Item(132, 13).Value = -2295.7142
$ws.Cells.Item(132, 14).Value = -449830.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 118772.12
$ws.Cells.Item(136, 9).Value = 91847.55
$ws.Cells.Item(136, 11).Value = 275542.65
$ws.Cells.Item(136, 13).Value = -272992.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2325.2307
$ws.Cells.Item(122, 9).Value = 1492.8823
$ws.Cells.Item(122, 10).Value = 2968.4092
$ws.Cells.Item(122, 11).Value = 4478.6469
$ws.Cells.Item(122, 12).Value = 8905.2276
$ws.Cells.Item(122, 13).Value = -2028.6469
$ws.Cells.Item(122, 14).Value = -13805.2276

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value = 46671.91
$ws.Cells.Item(123, 10).Value = 46671.91
$ws.Cells.Item(123, 12).Value = 46671.91
$ws.Cells.Item(123, 14).Value = -56471.91

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 75644.59
$ws.Cells.Item(132, 9).Value = 56647
$ws.Cells.Item(132, 10).Value = 113639.78
$ws.Cells.Item(132, 11).Value = 169941
$ws.Cells.Item(132, 12).Value = 340919.34
$ws.Cells.Item(132, 13).Value = -167411
$ws.Cells.Item(132, 14).Value = -345979.34

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 62824.637
$ws.Cells.Item(136, 9).Value = 38386.223
$ws.Cells.Item(136, 10).Value = 172797.5
$ws.Cells.Item(136, 11).Value = 115158.669
$ws.Cells.Item(136, 12).Value = 518392.5
$ws.Cells.Item(136, 13).Value = -112608.669
$ws.Cells.Item(136, 14).Value = -523492.5
